$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = '[0.16820465324541445, 13.613295181444823]'
$ws.Range("N2").Value = 0.04476261793504355
$ws.Range("O2").Value = 0.04476261793504355
$ws.Range("Q2").Value = '[-2.943474197958004, -0.2515789912784623]'
$ws.Range("R2").Value = 0.02106787524876497
$ws.Range("S2").Value = 0.02106787524876497
$ws.Range("U2").Value = '[5.179046595938832, 13.201658799610318]'
$ws.Range("V2").Value = 0.00003280012332829152
$ws.Range("W2").Value = 0.00003280012332829152
$ws.Range("Y2").Value = 1.035035035035062
$ws.Range("Z2").Value = 12.10990990991019

$ws.Range("M3").Value = '[-0.3520125859853618, 13.710080739939412]'
$ws.Range("N3").Value = 0.06208910788380839
$ws.Range("O3").Value = 0.06208910788380839
$ws.Range("Q3").Value = '[-3.144737390980774, -0.2515789912784623]'
$ws.Range("R3").Value = 0.02243492632858546
$ws.Range("S3").Value = 0.02243492632858546
$ws.Range("U3").Value = '[5.072917063133222, 13.054029384043528]'
$ws.Range("V3").Value = 0.00003737269203885063
$ws.Range("W3").Value = 0.00003737269203885063
$ws.Range("Y3").Value = 1.035035035035062
$ws.Range("Z3").Value = 12.93793793793824

$ws.Range("M4").Value = '[-1.7818358605128033, 15.093653572821847]'
$ws.Range("N4").Value = 0.1191131247534338
$ws.Range("O4").Value = 0.1191131247534338
$ws.Range("Q4").Value = '[-4.207658629132275, 0.05660527303765317]'
$ws.Range("R4").Value = 0.05612975575749246
$ws.Range("S4").Value = 0.05612975575749246
$ws.Range("U4").Value = '[5.6556041125666665, 14.41100855071905]'
$ws.Range("V4").Value = 0.00003262826343863878
$ws.Range("W4").Value = 0.00003262826343863878
$ws.Range("Y4").Value = -0.2328828828828851
$ws.Range("Z4").Value = 17.31096096096136

$ws.Range("M5").Value = '[-1.1135596456013062, 14.656233152894258]'
$ws.Range("N5").Value = 0.09054614305422826
$ws.Range("O5").Value = 0.09054614305422826
$ws.Range("Q5").Value = '[-5.346053564667316, 0.4151053356094625]'
$ws.Range("R5").Value = 0.09159769800753392
$ws.Range("S5").Value = 0.09159769800753392
$ws.Range("U5").Value = '[4.7925607956322605, 12.831211035966726]'
$ws.Range("V5").Value = 0.00006250185303646738
$ws.Range("W5").Value = 0.00006250185303646738
$ws.Range("Y5").Value = -1.707807807807848
$ws.Range("Z5").Value = 21.99449449449501

$ws.Range("M6").Value = '[0.08931264229155111, 14.110173810641937]'
$ws.Range("N6").Value = 0.04727036835786658
$ws.Range("O6").Value = 0.04727036835786658
$ws.Range("Q6").Value = '[-4.364895498681314, -1.3711055024676182]'
$ws.Range("R6").Value = 0.0003606969223073353
$ws.Range("S6").Value = 0.0003606969223073353
$ws.Range("U6").Value = '[5.188586629174342, 13.14887983136871]'
$ws.Range("V6").Value = 0.00003021033159988917
$ws.Range("W6").Value = 0.00003021033159988917
$ws.Range("Y6").Value = 4.971011011011043
$ws.Range("Z6").Value = 15.82514514514523

$ws.Range("B7").Value = 0
$ws.Range("M7").Value = '[-0.31174202527005335, 14.732967950099242]'
$ws.Range("N7").Value = 0.05984380628795671
$ws.Range("O7").Value = 0.05984380628795671
$ws.Range("Q7").Value = '[-4.302000750861699, -1.2830528555201566]'
$ws.Range("R7").Value = 0.0005406574194140035
$ws.Range("S7").Value = 0.0005406574194140035
$ws.Range("U7").Value = '[5.196254053252073, 13.12791661113242]'
$ws.Range("V7").Value = 0.00002891804873694781
$ws.Range("W7").Value = 0.00002891804873694781
$ws.Range("Y7").Value = 4.6517717717718
$ws.Range("Z7").Value = 15.5971171171172

$ws.Range("M8").Value = '[-0.21137966618823256, 14.618027891891309]'
$ws.Range("N8").Value = 0.05660679209161268
$ws.Range("O8").Value = 0.05660679209161268
$ws.Range("Q8").Value = '[-3.962369112635776, -0.3144737390980774]'
$ws.Range("R8").Value = 0.02259756870124141
$ws.Range("S8").Value = 0.02259756870124141
$ws.Range("U8").Value = '[5.009910123778655, 12.983219239702407]'
$ws.Range("V8").Value = 0.00004111237131176892
$ws.Range("W8").Value = 0.00004111237131176892
$ws.Range("Y8").Value = 1.140140140140145
$ws.Range("Z8").Value = 14.36576576576585

$ws.Range("M9").Value = '[-0.45386810048483284, 14.58806247459134]'
$ws.Range("N9").Value = 0.06485994897373826
$ws.Range("O9").Value = 0.06485994897373826
$ws.Range("Q9").Value = '[-3.89947436481616, -0.22642109215061534]'
$ws.Range("R9").Value = 0.02854745622337251
$ws.Range("S9").Value = 0.02854745622337251
$ws.Range("U9").Value = '[4.965439649399997, 12.973224348357078]'
$ws.Range("V9").Value = 0.00004579960710593411
$ws.Range("W9").Value = 0.00004579960710593411
$ws.Range("Y9").Value = 0.8209009009009032
$ws.Range("Z9").Value = 14.13773773773782

$ws.Range("M10").Value = '[0.0924782028511082, 15.080786589751286]'
$ws.Range("N10").Value = 0.04735404460824433
$ws.Range("O10").Value = 0.04735404460824433
$ws.Range("Q10").Value = '[-2.666737307551696, 0.2264210921506149]'
$ws.Range("R10").Value = 0.09625166952444952
$ws.Range("S10").Value = 0.09625166952444952
$ws.Range("U10").Value = '[5.323763774923859, 13.630777576002465]'
$ws.Range("V10").Value = 0.00003488018817887628
$ws.Range("W10").Value = 0.00003488018817887628
$ws.Range("Y10").Value = -0.8209009009009005
$ws.Range("Z10").Value = 9.668388388388442
